# "update scripts wuth new tpm"
# Re-run of the NATMI ligand/receptor scoring pipeline (Fgf22-Fgfr1) against
# updated TPM input: every derived-specificity / weighted-expression metric
# in columns I,J and M:T is refreshed with the new numbers. Sending/ligand/
# receptor/target-cluster labels (columns A-D) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.3046368955123587
$ws.Range("J2").Value = 0.3046368955123587
$ws.Range("M2").Value = 10.402079
$ws.Range("N2").Value = 20.804158
$ws.Range("O2").Value = 0.117441350183963
$ws.Range("P2").Value = 0.08862141909929068
$ws.Range("Q2").Value = 0.3269616144876667
$ws.Range("R2").Value = 1.961769686926
$ws.Range("S2").Value = 0.03577696832482225
$ws.Range("T2").Value = 0.02699735399030757
$ws.Range("I3").Value = 0.3046368955123587
$ws.Range("J3").Value = 0.3046368955123587
$ws.Range("O3").Value = 0.6341369869521791
$ws.Range("P3").Value = 0.7177810830557603
$ws.Range("S3").Value = 0.193181523034673
$ws.Range("T3").Value = 0.2186626007996053
$ws.Range("I4").Value = 0.3046368955123587
$ws.Range("J4").Value = 0.3046368955123587
$ws.Range("M4").Value = 0.2909543333333333
$ws.Range("N4").Value = 0.8728629999999999
$ws.Range("O4").Value = 0.003284926960133785
$ws.Range("P4").Value = 0.003718216220971988
$ws.Range("Q4").Value = 0.009145373590111111
$ws.Range("R4").Value = 0.082308362311
$ws.Range("S4").Value = 0.001000709951120006
$ws.Range("T4").Value = 0.001132705846400601
$ws.Range("I5").Value = 0.3046368955123587
$ws.Range("J5").Value = 0.3046368955123587
$ws.Range("M5").Value = 20.562391
$ws.Range("N5").Value = 41.124782
$ws.Range("O5").Value = 0.232153107282743
$ws.Range("P5").Value = 0.175183083160057
$ws.Range("Q5").Value = 0.6463239280423333
$ws.Range("R5").Value = 3.877943568254
$ws.Range("S5").Value = 0.07072240188616237
$ws.Range("T5").Value = 0.05336723060016313
$ws.Range("I6").Value = 0.3046368955123587
$ws.Range("J6").Value = 0.3046368955123587
$ws.Range("M6").Value = 0.5741476666666667
$ws.Range("N6").Value = 1.722443
$ws.Range("O6").Value = 0.006482230828885768
$ws.Range("P6").Value = 0.007337251667557973
$ws.Range("Q6").Value = 0.01804680084122222
$ws.Range("R6").Value = 0.162421207571
$ws.Range("S6").Value = 0.001974726675706264
$ws.Range("T6").Value = 0.002235197569597738
$ws.Range("I7").Value = 0.3046368955123587
$ws.Range("J7").Value = 0.3046368955123587
$ws.Range("M7").Value = 0.5758453333333333
$ws.Range("N7").Value = 1.727536
$ws.Range("O7").Value = 0.006501397792095299
$ws.Range("P7").Value = 0.00735894679636216
$ws.Range("Q7").Value = 0.01810016246577778
$ws.Range("R7").Value = 0.162901462192
$ws.Range("S7").Value = 0.001980565639874815
$ws.Range("T7").Value = 0.002241806706284386
$ws.Range("G8").Value = 0.07174733333333333
$ws.Range("H8").Value = 0.215242
$ws.Range("I8").Value = 0.6953631044876413
$ws.Range("J8").Value = 0.6953631044876413
$ws.Range("M8").Value = 10.402079
$ws.Range("N8").Value = 20.804158
$ws.Range("O8").Value = 0.117441350183963
$ws.Range("P8").Value = 0.08862141909929068
$ws.Range("Q8").Value = 0.7463214293726667
$ws.Range("R8").Value = 4.477928576236
$ws.Range("S8").Value = 0.08166438185914071
$ws.Range("T8").Value = 0.06162406510898311
$ws.Range("I9").Value = 0.6953631044876413
$ws.Range("J9").Value = 0.6953631044876413
$ws.Range("O9").Value = 0.6341369869521791
$ws.Range("P9").Value = 0.7177810830557603
$ws.Range("Q9").Value = 4.029841463665777
$ws.Range("R9").Value = 36.268573172992
$ws.Range("S9").Value = 0.4409554639175061
$ws.Range("T9").Value = 0.499118482256155
$ws.Range("I10").Value = 0.6953631044876413
$ws.Range("J10").Value = 0.6953631044876413
$ws.Range("M10").Value = 0.2909543333333333
$ws.Range("N10").Value = 0.8728629999999999
$ws.Range("O10").Value = 0.003284926960133785
$ws.Range("P10").Value = 0.003718216220971988
$ws.Range("Q10").Value = 0.02087519753844444
$ws.Range("R10").Value = 0.187876777846
$ws.Range("S10").Value = 0.002284217009013779
$ws.Range("T10").Value = 0.002585510374571387
$ws.Range("I11").Value = 0.6953631044876413
$ws.Range("J11").Value = 0.6953631044876413
$ws.Range("M11").Value = 20.562391
$ws.Range("N11").Value = 41.124782
$ws.Range("O11").Value = 0.232153107282743
$ws.Range("P11").Value = 0.175183083160057
$ws.Range("Q11").Value = 1.475296721207333
$ws.Range("R11").Value = 8.851780327243999
$ws.Range("S11").Value = 0.1614307053965806
$ws.Range("T11").Value = 0.1218158525598939
$ws.Range("I12").Value = 0.6953631044876413
$ws.Range("J12").Value = 0.6953631044876413
$ws.Range("M12").Value = 0.5741476666666667
$ws.Range("N12").Value = 1.722443
$ws.Range("O12").Value = 0.006482230828885768
$ws.Range("P12").Value = 0.007337251667557973
$ws.Range("Q12").Value = 0.04119356402288889
$ws.Range("R12").Value = 0.370742076206
$ws.Range("S12").Value = 0.004507504153179505
$ws.Range("T12").Value = 0.005102054097960235
$ws.Range("I13").Value = 0.6953631044876413
$ws.Range("J13").Value = 0.6953631044876413
$ws.Range("M13").Value = 0.5758453333333333
$ws.Range("N13").Value = 1.727536
$ws.Range("O13").Value = 0.006501397792095299
$ws.Range("P13").Value = 0.00735894679636216
$ws.Range("Q13").Value = 0.04131536707911111
$ws.Range("R13").Value = 0.371838303712
$ws.Range("S13").Value = 0.004520832152220483
$ws.Range("T13").Value = 0.005117140090077773